$wb = $excel.ActiveWorkbook

# --- safety_orders: remove the first safety-order row (old row 2, Safety Order No. 3) ---
# Remaining rows shift up, matching rows 3-6 -> 2-5.
$wsSafety = $wb.Worksheets.Item("safety_orders")
$wsSafety.Rows.Item(2).Delete()

# --- open_buy_orders: append a new available-order row ---
$wsBuy = $wb.Worksheets.Item("open_buy_orders")
$wsBuy.Cells.Item(4, 1).Value = "OGGP6M-BDBT2-XIB2QR"
$wsBuy.Cells.Item(4, 2).Value = 1.62216
